# Performance.xlsx: add "resilient" and "manhattan" sheets (parts of the
# Part 3 experiment results), matching the "extended"/"original" sheet
# layout, and update sheet selections/active tab.

$wb = $excel.ActiveWorkbook

$wsOriginal = $wb.Worksheets.Item(1)
$wsExtended = $wb.Worksheets.Item(2)

# --- update the "extended" sheet's current selection -----------------
[void]$wsExtended.Activate()
$wsExtended.Range("A2:A4").Select() | Out-Null

# --- add the two new worksheets, appended after "extended" -----------
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "resilient"

$ws4 = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "manhattan"

# --- column A width / labels ------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 26.166666666666668
$ws4.Columns.Item(1).ColumnWidth = 25.761904761904763

$ws3.Range("A2").Value = "Epochs"
$ws3.Range("A3").Value = "Error"
$ws3.Range("A4").Value = "Correctly classified characters"
$ws3.Range("L1").Value = "avg"

$ws4.Range("A2").Value = "Epochs"
$ws4.Range("A3").Value = "Error"
$ws4.Range("A4").Value = "Correctly classified characters"
$ws4.Range("L1").Value = "avg"

# --- numeric data -------------------------------------------------------
$ws3.Cells.Item(2, 2).Value = 17
$ws3.Cells.Item(2, 3).Value = 14
$ws3.Cells.Item(2, 4).Value = 18
$ws3.Cells.Item(2, 5).Value = 18
$ws3.Cells.Item(2, 6).Value = 15
$ws3.Cells.Item(2, 7).Value = 15
$ws3.Cells.Item(2, 8).Value = 16
$ws3.Cells.Item(2, 9).Value = 15
$ws3.Cells.Item(2, 10).Value = 16
$ws3.Cells.Item(2, 11).Value = 15
$ws3.Cells.Item(3, 2).Value = 0.00765535887213095
$ws3.Cells.Item(3, 3).Value = 0.00820308014416701
$ws3.Cells.Item(3, 4).Value = 0.0095119891508798
$ws3.Cells.Item(3, 5).Value = 0.00553087859751923
$ws3.Cells.Item(3, 6).Value = 0.00603335839697162
$ws3.Cells.Item(3, 7).Value = 0.00799146288550538
$ws3.Cells.Item(3, 8).Value = 0.00682427245060408
$ws3.Cells.Item(3, 9).Value = 0.00849409135225562
$ws3.Cells.Item(3, 10).Value = 0.00819632954536385
$ws3.Cells.Item(3, 11).Value = 0.00989665898469715
$ws3.Cells.Item(4, 2).Value = 100
$ws3.Cells.Item(4, 3).Value = 100
$ws3.Cells.Item(4, 4).Value = 100
$ws3.Cells.Item(4, 5).Value = 100
$ws3.Cells.Item(4, 6).Value = 100
$ws3.Cells.Item(4, 7).Value = 100
$ws3.Cells.Item(4, 8).Value = 100
$ws3.Cells.Item(4, 9).Value = 100
$ws3.Cells.Item(4, 10).Value = 100
$ws3.Cells.Item(4, 11).Value = 100
$ws4.Cells.Item(2, 2).Value = 114
$ws4.Cells.Item(2, 3).Value = 57
$ws4.Cells.Item(2, 4).Value = 35
$ws4.Cells.Item(2, 5).Value = 40
$ws4.Cells.Item(2, 6).Value = 88
$ws4.Cells.Item(2, 7).Value = 1229
$ws4.Cells.Item(2, 8).Value = 54
$ws4.Cells.Item(2, 9).Value = 51
$ws4.Cells.Item(2, 10).Value = 33
$ws4.Cells.Item(2, 11).Value = 34
$ws4.Cells.Item(3, 2).Value = 0.00918759642631555
$ws4.Cells.Item(3, 3).Value = 0.00925018759132901
$ws4.Cells.Item(3, 4).Value = 0.00821787300942242
$ws4.Cells.Item(3, 5).Value = 0.00320922788091682
$ws4.Cells.Item(3, 6).Value = 0.00758393346269152
$ws4.Cells.Item(3, 7).Value = 0.00995156067338492
$ws4.Cells.Item(3, 8).Value = 0.00905018657787189
$ws4.Cells.Item(3, 9).Value = 0.00911626184025707
$ws4.Cells.Item(3, 10).Value = 0.0089289267468091
$ws4.Cells.Item(3, 11).Value = 0.00809009476906786
$ws4.Cells.Item(4, 2).Value = 95.2380952380952
$ws4.Cells.Item(4, 3).Value = 80.9523809523809
$ws4.Cells.Item(4, 4).Value = 92.8571428571428
$ws4.Cells.Item(4, 5).Value = 64.2857142857142
$ws4.Cells.Item(4, 6).Value = 100
$ws4.Cells.Item(4, 7).Value = 85.7142857142857
$ws4.Cells.Item(4, 8).Value = 97.6190476190476
$ws4.Cells.Item(4, 9).Value = 100
$ws4.Cells.Item(4, 10).Value = 97.6190476190476
$ws4.Cells.Item(4, 11).Value = 92.8571428571428

# --- average formulas ----------------------------------------------------
$ws3.Cells.Item(2, 12).Formula = "=AVERAGE(B2:K2)"
$ws3.Cells.Item(3, 12).Formula = "=AVERAGE(B3:K3)"
$ws3.Cells.Item(4, 12).Formula = "=AVERAGE(B4:K4)"

$ws4.Cells.Item(2, 12).Formula = "=AVERAGE(B2:K2)"
$ws4.Cells.Item(3, 12).Formula = "=AVERAGE(B3:K3)"
$ws4.Cells.Item(4, 12).Formula = "=AVERAGE(B4:K4)"

# --- sheet view / selection state ----------------------------------------
[void]$ws3.Activate()
$ws3.Range("K4").Select() | Out-Null

[void]$ws4.Activate()
$ws4.Range("K5").Select() | Out-Null

Write-Host "done"
